$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.71"
$ws.Range("E2").Value = "'0.36%"
$ws.Range("G2").Value = "'23"

$ws.Range("D3").Value = "'35.61"
$ws.Range("E3").Value = "'10.99%"
$ws.Range("G3").Value = "'23"

$ws.Range("D4").Value = "'5.088"
$ws.Range("E4").Value = "'1.40%"
$ws.Range("G4").Value = "'23"

$ws.Range("D5").Value = "'0.07757"
$ws.Range("E5").Value = "'-0.66%"
$ws.Range("G5").Value = "'23"

$ws.Range("D6").Value = "'2.216"
$ws.Range("E6").Value = "'-7.59%"
$ws.Range("G6").Value = "'23"

$ws.Range("D7").Value = "'8.040"
$ws.Range("E7").Value = "'0.52%"
$ws.Range("G7").Value = "'23"

$ws.Range("E8").Value = "'4.16%"
$ws.Range("G8").Value = "'23"

$ws.Range("D9").Value = "'0.9276"
$ws.Range("E9").Value = "'-0.90%"
$ws.Range("G9").Value = "'23"

$ws.Range("D10").Value = "'0.09566"
$ws.Range("E10").Value = "'-5.98%"
$ws.Range("G10").Value = "'23"

$ws.Range("E11").Value = "'2.89%"
$ws.Range("G11").Value = "'23"

$ws.Range("D12").Value = "'0.08516"
$ws.Range("E12").Value = "'0.19%"
$ws.Range("G12").Value = "'23"

$ws.Range("D13").Value = "'0.03642"
$ws.Range("E13").Value = "'9.12%"
$ws.Range("G13").Value = "'23"

$ws.Range("D14").Value = "'0.09917"
$ws.Range("E14").Value = "'-0.02%"
$ws.Range("G14").Value = "'23"

$ws.Range("D15").Value = "'0.001479"
$ws.Range("E15").Value = "'-0.68%"
$ws.Range("G15").Value = "'23"

$ws.Range("D16").Value = "'0.005743"
$ws.Range("E16").Value = "'-0.47%"
$ws.Range("G16").Value = "'23"

$ws.Range("D17").Value = "'3.487"
$ws.Range("E17").Value = "'-0.29%"
$ws.Range("G17").Value = "'23"

$ws.Range("D18").Value = "'2.183"
$ws.Range("E18").Value = "'0.03%"
$ws.Range("G18").Value = "'23"

$ws.Range("E19").Value = "'2.93%"
$ws.Range("G19").Value = "'23"

$ws.Range("D20").Value = "'0.1323"
$ws.Range("E20").Value = "'-1.51%"
$ws.Range("G20").Value = "'23"

$ws.Range("D21").Value = "'4.562"
$ws.Range("E21").Value = "'5.68%"
$ws.Range("G21").Value = "'23"

$ws.Range("D22").Value = "'0.2243"
$ws.Range("E22").Value = "'7.32%"
$ws.Range("G22").Value = "'23"

$ws.Range("D23").Value = "'0.04681"
$ws.Range("E23").Value = "'1.61%"
$ws.Range("G23").Value = "'23"

$ws.Range("D24").Value = "'0.001236"
$ws.Range("E24").Value = "'1.42%"
$ws.Range("G24").Value = "'23"

$ws.Range("D25").Value = "'0.004508"
$ws.Range("E25").Value = "'2.55%"
$ws.Range("G25").Value = "'23"

$ws.Range("D26").Value = "'0.0001306"
$ws.Range("E26").Value = "'0.84%"
$ws.Range("G26").Value = "'23"

$ws.Range("E27").Value = "'-20.07%"
$ws.Range("G27").Value = "'23"

$ws.Range("G28").Value = "'23"

$ws.Range("G29").Value = "'23"

$ws.Range("G30").Value = "'23"

$ws.Range("G31").Value = "'23"

$ws.Range("G32").Value = "'23"

$ws.Range("G33").Value = "'23"

$ws.Range("G34").Value = "'23"

$ws.Range("G35").Value = "'23"

$ws.Range("G36").Value = "'23"

$ws.Range("G37").Value = "'23"

$ws.Range("G38").Value = "'23"

$ws.Range("D39").Value = "'0.01758"
$ws.Range("E39").Value = "'2.07%"
$ws.Range("G39").Value = "'23"

$ws.Range("D40").Value = "'0.04729"
$ws.Range("E40").Value = "'-1.98%"
$ws.Range("G40").Value = "'23"

$ws.Range("D41").Value = "'0.007954"
$ws.Range("E41").Value = "'2.66%"
$ws.Range("G41").Value = "'23"

$ws.Range("D42").Value = "'0.1408"
$ws.Range("E42").Value = "'0.20%"
$ws.Range("G42").Value = "'23"

$ws.Range("D43").Value = "'0.007878"
$ws.Range("E43").Value = "'-19.50%"
$ws.Range("G43").Value = "'23"

$ws.Range("D44").Value = "'0.002231"
$ws.Range("E44").Value = "'5.65%"
$ws.Range("G44").Value = "'23"

$ws.Range("D45").Value = "'0.009646"
$ws.Range("E45").Value = "'-4.89%"
$ws.Range("G45").Value = "'23"

$ws.Range("D46").Value = "'0.00006209"
$ws.Range("E46").Value = "'1.82%"
$ws.Range("G46").Value = "'23"

$ws.Range("D47").Value = "'0.00000000753"
$ws.Range("E47").Value = "'0.92%"
$ws.Range("G47").Value = "'23"

$ws.Range("D48").Value = "'5.787"
$ws.Range("E48").Value = "'117.98%"
$ws.Range("G48").Value = "'23"

$ws.Range("D49").Value = "'0.002702"
$ws.Range("E49").Value = "'35.74%"
$ws.Range("G49").Value = "'23"

$ws.Range("D50").Value = "'0.00002109"
$ws.Range("E50").Value = "'0.92%"
$ws.Range("G50").Value = "'23"

$ws.Range("D51").Value = "'0.0002009"
$ws.Range("E51").Value = "'0.92%"
$ws.Range("G51").Value = "'23"
